$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for several rows per the repull/recalculated data
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = -3
$ws.Range("F11").Value = -2
$ws.Range("F13").Value = -5
